# Add a "Spain" test-data sheet, cloned from the existing "Italy" sheet
# (same layout/styles), with the market name and ticket reference swapped
# in, positioned right after "Italy".

$wb = $excel.ActiveWorkbook

$italy = $wb.Worksheets.Item("Italy")

# Clone Italy -> new sheet lands immediately after Italy.
$italy.Copy($null, $italy)
$spain = $wb.Worksheets.Item($wb.Worksheets.Count)
$spain.Name = "Spain"

# Swap in the Spain-specific values.
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3103/T2034/T2037 "

# Widen columns A, B, D to fit the new content (best-fit-style autosize).
$spain.Columns.Item(1).ColumnWidth = 24.27
$spain.Columns.Item(2).ColumnWidth = 33.62
$spain.Columns.Item(4).ColumnWidth = 16.28

# Taller rows 3-5 to accommodate the wrapped text.
$spain.Range("A3:A5").EntireRow.RowHeight = 28.8

# Move the active selection off Italy (deselect its tab) and onto the
# new Spain sheet, which becomes the active tab.
$italy.Range("A1:D11").Select() | Out-Null
$spain.Range("D14").Select() | Out-Null
